# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches from the custom "Table_0" style to the
#     built-in table style {5747DB66-5EBB-4A26-B32B-22DB7367801D}.
#  2. The presentation's theme (ppt/theme/theme2.xml, linked from the
#     slide master) is recoloured from the "Integral" palette to the
#     stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$s6 = $p.Slides.Item(6)
$tbl = $s6.Shapes.Item(2).Table
$tbl.ApplyStyle("{5747DB66-5EBB-4A26-B32B-22DB7367801D}")

# --- 2. Recolour the deck's theme to the "Office Theme" palette ---------
$theme = $p.SlideMaster.Theme
$clrScheme = $theme.ThemeColorScheme

# Order matches the ThemeColorScheme index contract:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $clrScheme.Item($i).RGB = $officeThemeRGB[$i - 1]
}
